$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 0
$ws.Range("A41").Value = 0
$ws.Range("A43").Value = 0
$ws.Range("A54").Value = 0
$ws.Range("A55").Value = 0
$ws.Range("A56").Value = 0
$ws.Range("A308").Value = 0
